$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells I1 ("I0") and J1 ("IF") ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, centered, bordered) from H1 onto I1 and J1
# so the new header cells reuse the existing cell format exactly.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows 2-33: I column = 1, J column = copy of H column value ---
for ($r = 2; $r -le 33; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value()
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
